$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.314.51"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.866.51"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'235.63"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4673"
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("D8").Value = "'0.2839"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "'0.06521"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'21.90"
$ws.Range("E10").Value = "  +9.03%  "
$ws.Range("D11").Value = "'0.07936"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").Value = "'97.31"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "1.868.90"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "'5.155"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "'0.6779"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "'279.79"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "30.312.28"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'13.27"
$ws.Range("E18").Value = "  +5.24%  "
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'5.418"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Value = "2.116.39"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").Value = "'0.000007310"
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'6.153"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "'166.32"
$ws.Range("E25").Value = "  -0.94%  "
$ws.Range("D26").Value = "'9.158"
$ws.Range("E26").Value = "  -1.65%  "
$ws.Range("D27").Value = "'19.07"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'1.933"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'1.388"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("D30").Value = "'0.09725"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "'4.400"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'1.477"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").Value = "'4.086"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'0.04748"
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").Value = "'1.133"
$ws.Range("E35").Value = "  +4.13%  "
$ws.Range("D36").Value = "'0.7076"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.01867"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "'2.573"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "'6.322"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'74.66"
$ws.Range("E41").Value = "  +3.70%  "
$ws.Range("D42").Value = "'1.970"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").Value = "'0.8505"
$ws.Range("E43").Value = "  -1.29%  "
$ws.Range("D44").Value = "'0.4189"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("D46").Value = "'103.33"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").Value = "'970.94"
$ws.Range("E47").Value = "  -5.38%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.193"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.336"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "'34.11"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").Value = "'0.1132"
$ws.Range("E51").Value = "  -0.98%  "
